$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet to match the new ASIN ---
$ws.Name = "B08K4353V1"

# --- Replace the 100 keyword values in column A (A1:A100) ---
$values = @(
    "neck compression wrap",
    "training bra sets",
    "wide waistband yoga",
    "girls seamless underwear",
    "sport bra",
    "seamless bra",
    "yoga legging",
    "workout waist",
    "yoga workout",
    "yoga bra",
    "sports bra for women",
    "sport apparel",
    "seamless bra set",
    "workout apparel",
    "workout bra",
    "workout outfit set",
    "workout set",
    "workout leggings for women high waist",
    "workout legging",
    "legging bra set",
    "workout outfit",
    "sports leggings women",
    "black legging",
    "sport bra set",
    "yoga sets for women",
    "high waist workout",
    "yoga set purple",
    "yoga bras for women",
    "2 pcs",
    "green 2",
    "pcs set",
    "bra set",
    "waist yoga",
    "outfit set",
    "seamless sports bra",
    "seamless workout leggings",
    "seamless yoga bra",
    "seamless yoga",
    "sports for women",
    "high waist yoga",
    "sport legging",
    "black bra",
    "yoga legging set",
    "yoga workout set",
    "yoga bra set",
    "yoga outfit set",
    "workout set women",
    "black legging set",
    "sports leggings for women",
    "black bra set",
    "sport leggings for women",
    "sport set",
    "black outfit",
    "high waist",
    "set 2",
    "yoga apparel",
    "purple bra",
    "purple legging",
    "high waist legging",
    "high sport bra",
    "black n",
    "womens sports bra set",
    "sport outfit",
    "sport waist",
    "sport workout",
    "waist bra",
    "black 2",
    "yoga sport",
    "n set",
    "green set",
    "yoga 2",
    "green bra",
    "green bra set",
    "green yoga",
    "yoga sport bra",
    "womens workout sports bra",
    "high waist set",
    "ready set",
    "high black",
    "waist set",
    "black yoga",
    "green outfit",
    "black apparel",
    "purple yoga",
    "black set",
    "purple outfit",
    "purple workout",
    "black sport",
    "seamless legging",
    "green apparel",
    "purple set",
    "black seamless bra",
    "black yoga bra",
    "seamless workout",
    "green sport bra",
    "black sport bra",
    "purple sport bra",
    "workout sport bra",
    "black workout bra",
    "seamless sport bra"
)

$arr = New-Object 'object[,]' 100,1
for ($i = 0; $i -lt 100; $i++) {
    $arr[$i, 0] = $values[$i]
}
$ws.Range("A1:A100").Value = $arr

# --- Clear the explicit left alignment on the data style, reverting to General ---
$ws.Range("A1:A100").HorizontalAlignment = 1

# --- Re-create (then discard) the "highlight duplicates" rule twice; Excel keeps
#     the now-orphaned style records (dxfs) around even after the rule is removed,
#     which is why dxfs count grows from 1 to 3 while the sheet still shows a
#     single duplicate-values conditional format. ---
for ($i = 0; $i -lt 2; $i++) {
    $tmp = $ws.Range("A1:A100").FormatConditions.AddUniqueValues()
    $tmp.DupeUnique = 1
    $tmp.Interior.Color = 13998939
    $tmp.Delete()
}

# --- Update the view: scroll so row 10 is at the top and select the full range ---
$ws.Range("A1:A100").Select()
$excel.ActiveWindow.ScrollRow = 10
